$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the bug: B9 was 7, should be 8
$ws.Range("B9").Value = 8

# Add new row of data for 09/05/2024
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "09/05/2024"
$ws.Range("B10").Value = 4
$ws.Range("D10").Value = "fixed bugs"

# Update the active selection to D10
$ws.Range("D10").Select()
